$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style/format from row 301 template down to new rows 302:328
$ws.Range("A301:D301").Copy($ws.Range("A302:D328"))

$data = @(
  @(302, 44376, 0, 4, 15.6561900661474),
  @(303, 44377, 1, 5, 19.57023758268425),
  @(304, 44378, 1, 5, 19.57023758268425),
  @(305, 44379, 0, 5, 19.57023758268425),
  @(306, 44380, 0, 5, 19.57023758268425),
  @(307, 44381, 0, 3, 11.74214254961055),
  @(308, 44382, 1, 3, 11.74214254961055),
  @(309, 44383, 0, 3, 11.74214254961055),
  @(310, 44384, 0, 2, 7.828095033073701),
  @(311, 44385, 2, 3, 11.74214254961055),
  @(312, 44386, 0, 3, 11.74214254961055),
  @(313, 44387, 1, 4, 15.6561900661474),
  @(314, 44388, 0, 4, 15.6561900661474),
  @(315, 44389, 0, 3, 11.74214254961055),
  @(316, 44390, 0, 3, 11.74214254961055),
  @(317, 44391, 0, 3, 11.74214254961055),
  @(318, 44392, 0, 1, 3.914047516536851),
  @(319, 44393, 0, 1, 3.914047516536851),
  @(320, 44394, 0, 0, 0),
  @(321, 44395, 0, 0, 0),
  @(322, 44396, 0, 0, 0),
  @(323, 44397, 0, 0, 0),
  @(324, 44398, 0, 0, 0),
  @(325, 44399, 0, 0, 0),
  @(326, 44400, 0, 0, 0),
  @(327, 44401, 0, 0, 0),
  @(328, 44402, 1, 1, 3.914047516536851)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "done"
